# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# --- Update the re-generated query timestamps on the "data" sheet (column F) ---
$timestamps = @(
    "2021-10-05 14:21:34.079104"
    "2021-10-05 14:21:34.079112"
    "2021-10-05 14:21:34.079115"
    "2021-10-05 14:21:34.079118"
    "2021-10-05 14:21:34.079121"
    "2021-10-05 14:21:34.079124"
    "2021-10-05 14:21:34.079127"
    "2021-10-05 14:21:34.079129"
    "2021-10-05 14:21:34.079132"
    "2021-10-05 14:21:34.079135"
    "2021-10-05 14:21:34.079137"
    "2021-10-05 14:21:34.079140"
    "2021-10-05 14:21:34.079143"
    "2021-10-05 14:21:34.079146"
    "2021-10-05 14:21:34.079148"
    "2021-10-05 14:21:34.079151"
    "2021-10-05 14:21:34.079153"
    "2021-10-05 14:21:34.079156"
    "2021-10-05 14:21:34.079159"
    "2021-10-05 14:21:34.079161"
    "2021-10-05 14:21:34.079164"
    "2021-10-05 14:21:34.079166"
    "2021-10-05 14:21:34.079169"
    "2021-10-05 14:21:34.079171"
    "2021-10-05 14:21:34.079175"
    "2021-10-05 14:21:34.079177"
    "2021-10-05 14:21:34.079180"
    "2021-10-05 14:21:34.079183"
    "2021-10-05 14:21:34.079185"
    "2021-10-05 14:21:34.079188"
    "2021-10-05 14:21:34.079190"
    "2021-10-05 14:21:34.079193"
    "2021-10-05 14:21:34.079196"
    "2021-10-05 14:21:34.079199"
    "2021-10-05 14:21:34.079201"
    "2021-10-05 14:21:34.079204"
    "2021-10-05 14:21:34.079207"
    "2021-10-05 14:21:34.079209"
    "2021-10-05 14:21:34.079212"
    "2021-10-05 14:21:34.079215"
    "2021-10-05 14:21:34.079218"
)

for ($i = 0; $i -lt $timestamps.Count; $i++) {
    $row = $i + 2
    $dataWs.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- Add the new "metadata" sheet, positioned right after "data" ---
$meta = $wb.Worksheets.Add($null, $dataWs)
$meta.Name = "metadata"

# Header row (bold, centered, top-aligned, thin border) - matches the "data" sheet header style
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $cell = $meta.Cells.Item(1, $c + 2)   # starts at column B
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108     # xlCenter
    $cell.VerticalAlignment = -4160       # xlTop
    $cell.Borders.LineStyle = 1           # xlContinuous
}

# Data row 2
$meta.Cells.Item(2, 1).Value = 0                                     # A2 - index
$meta.Cells.Item(2, 1).HorizontalAlignment = -4108
$meta.Cells.Item(2, 1).VerticalAlignment = -4160
$meta.Cells.Item(2, 1).Borders.LineStyle = 1

$meta.Cells.Item(2, 2).Value = "Mitochondrial disorder with complex IV deficiency"   # B2 - data_name
$meta.Cells.Item(2, 3).Value = 537                                                  # C2 - data_id
$meta.Cells.Item(2, 4).Value = "'1.12"                                              # D2 - data_version (kept as text)
$meta.Cells.Item(2, 5).Value = "2021-06-07T11:24:32.192970Z"                        # E2 - data_version_created
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:21:34.075601"                         # F2 - panel_query_time
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/537/?format=json"  # G2 - panel_get_request

$meta.Range("A1").Select() | Out-Null
